# Rotate the betting-odds rows for the three Azerbaijan Premier League
# fixtures stored on rows 172-174 (excluding the id/Div/Date columns
# A, C and D which stay put). The row that used to be on line 172 moves
# down to line 174, the row that used to be on 173 moves up to 172, and
# the row that used to be on 174 moves up to 173 (a 3-way cyclic shift).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B through AD (everything except A/C/D which are unchanged)
$cols = @("B","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")

$rows = @(172,173,174)

# Snapshot current values for the affected columns/rows before we start
# overwriting them.
$snapshot = @{}
foreach ($r in $rows) {
    $snapshot[$r] = @{}
    foreach ($c in $cols) {
        $snapshot[$r][$c] = $ws.Range("$c$r").Value2
    }
}

# New row 172 <- old row 173
# New row 173 <- old row 174
# New row 174 <- old row 172
$mapping = @{
    172 = 173
    173 = 174
    174 = 172
}

foreach ($destRow in $rows) {
    $srcRow = $mapping[$destRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value2 = $snapshot[$srcRow][$c]
    }
}
